$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1. Update "descriptives" sheet (sheet1) values
# -----------------------------------------------------------------
$wsDesc = $wb.Worksheets.Item("descriptives")

$wsDesc.Range("C2").Value = 845
$wsDesc.Range("D2").Value = 0.008030863400722361
$wsDesc.Range("E2").Value = 0.05060100631521103
$wsDesc.Range("F2").Value = 97.73350421127832
$wsDesc.Range("G2").Value = 84.34685245008768
$wsDesc.Range("H2").Value = 13.38665176119064

$wsDesc.Range("C3").Value = 438
$wsDesc.Range("D3").Value = 0.003975373526203995
$wsDesc.Range("E3").Value = 0.05301294003157205
$wsDesc.Range("F3").Value = 95.10650554477499
$wsDesc.Range("G3").Value = 88.47209472071785
$wsDesc.Range("H3").Value = 6.634410824057147

# -----------------------------------------------------------------
# 2. Update "coefficients" sheet (sheet2) values
# -----------------------------------------------------------------
$wsCoef = $wb.Worksheets.Item("coefficients")

$wsCoef.Range("C2").Value = 8.179405004212979
$wsCoef.Range("D2").Value = 0.7649727984206065
$wsCoef.Range("F2").Value = 9.764479940363771
$wsCoef.Range("G2").Value = 0.003630986766515545

$wsCoef.Range("C3").Value = 0.1962197121711471
$wsCoef.Range("D3").Value = 0.8333369649535934
$wsCoef.Range("F3").Value = 10.00026148235636
$wsCoef.Range("G3").Value = 0.8965894059026196

# -----------------------------------------------------------------
# 3. Add new "nr_studies" sheet (sheet3) at the end
# -----------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNr = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsNr.Name = "nr_studies"

# Header row
$headers = @("outcome", "quality_score_out_of_5", "n_effect_sizes", "k_studies")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $wsNr.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
}

# Data rows
$data = @(
    @("NS", "3", 300, 31),
    @("NS", "4", 393, 31),
    @("NS", "2", 74, 8),
    @("NS", "5", 69, 9),
    @("NS", "1", 9, 3),
    @("NT", "3", 149, 18),
    @("NT", "4", 195, 21),
    @("NT", "2", 41, 4),
    @("NT", "5", 53, 8)
)

# Column B (quality_score_out_of_5) holds text-typed digits in the source
# workbook (inlineStr "3", "4", ...), not numbers - force text format so the
# values are not coerced to numeric.
$wsNr.Range("B2:B10").NumberFormat = "@"

$r = 2
foreach ($row in $data) {
    $wsNr.Cells.Item($r, 1).Value = $row[0]
    $wsNr.Cells.Item($r, 2).Value = $row[1]
    $wsNr.Cells.Item($r, 3).Value = $row[2]
    $wsNr.Cells.Item($r, 4).Value = $row[3]
    $r++
}
